$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 7.084288597106934
$ws.Range("B1").Value = 5.274005889892578
$ws.Range("C1").Value = 4.436309814453125
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 3.413908004760742
